$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.4935314059257507
$ws.Range("B1").Value = 1.227790951728821
$ws.Range("C1").Value = 4.366623878479004
$ws.Range("D1").Value = 1.323935866355896
$ws.Range("E1").Value = 1.169966578483582
